$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Header row: add "Notes" to E1 ---
$ws.Range("E1").Value = "Notes"

# --- Remove underline rich-text formatting from D1 ("phase-in") ---
$ws.Range("D1").Characters(14, 8).Font.Underline = $false

# --- Rows 2-14: overwrite in place (keep existing style), new quarters 2013Q1-2016Q1 ---
$ws.Range("A2").Value = "2013Q1"
$ws.Range("D2").Value = 12.5
$ws.Range("E2").Value = "Pro forma, post-dividend"
$ws.Range("A3").Value = "2013Q2"
$ws.Range("D3").Value = 13.7
$ws.Range("A4").Value = "2013Q3"
$ws.Range("D4").Value = 13.5
$ws.Range("A5").Value = "2013Q4"
$ws.Range("D5").Value = 14
$ws.Range("A6").Value = "2014Q1"
$ws.Range("D6").Value = 10.7
$ws.Range("A7").Value = "2014Q2"
$ws.Range("D7").Value = 11.1
$ws.Range("A8").Value = "2014Q3"
$ws.Range("D8").Value = 12
$ws.Range("A9").Value = "2014Q4"
$ws.Range("D9").Value = 13
$ws.Range("A10").Value = "2015Q1"
$ws.Range("D10").Value = 13.4
$ws.Range("A11").Value = "2015Q2"
$ws.Range("D11").Value = 13.3
$ws.Range("A12").Value = "2015Q3"
$ws.Range("D12").Value = 13.7
$ws.Range("A13").Value = "2015Q4"
$ws.Range("D13").Value = 13
$ws.Range("A14").Value = "2016Q1"
$ws.Range("D14").Value = 12.8

# --- Rows 15-33: new unstyled rows, quarters 2016Q2-2020Q4 ---
$ws.Range("A15").Value = "2016Q2"
$ws.Range("B15").Value = "LYG"
$ws.Range("C15").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D15").Value = 13
$ws.Range("A16").Value = "2016Q3"
$ws.Range("B16").Value = "LYG"
$ws.Range("C16").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D16").Value = 14.1
$ws.Range("A17").Value = "2016Q4"
$ws.Range("B17").Value = "LYG"
$ws.Range("C17").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D17").Value = 13.8
$ws.Range("A18").Value = "2017Q1"
$ws.Range("B18").Value = "LYG"
$ws.Range("C18").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D18").Value = 14.5
$ws.Range("A19").Value = "2017Q2"
$ws.Range("B19").Value = "LYG"
$ws.Range("C19").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D19").Value = 14
$ws.Range("A20").Value = "2017Q3"
$ws.Range("B20").Value = "LYG"
$ws.Range("C20").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D20").Value = 14.9
$ws.Range("A21").Value = "2017Q4"
$ws.Range("B21").Value = "LYG"
$ws.Range("C21").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D21").Value = 15.5
$ws.Range("A22").Value = "2018Q1"
$ws.Range("B22").Value = "LYG"
$ws.Range("C22").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D22").Value = 16
$ws.Range("E22").Value = "Implied from prev quarter"
$ws.Range("A23").Value = "2018Q2"
$ws.Range("B23").Value = "LYG"
$ws.Range("C23").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D23").Value = 15.1
$ws.Range("A24").Value = "2018Q3"
$ws.Range("B24").Value = "LYG"
$ws.Range("C24").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D24").Value = 14.6
$ws.Range("A25").Value = "2018Q4"
$ws.Range("B25").Value = "LYG"
$ws.Range("C25").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D25").Value = 13.9
$ws.Range("A26").Value = "2019Q1"
$ws.Range("B26").Value = "LYG"
$ws.Range("C26").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D26").Value = 13.9
$ws.Range("A27").Value = "2019Q2"
$ws.Range("B27").Value = "LYG"
$ws.Range("C27").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D27").Value = 14
$ws.Range("A28").Value = "2019Q3"
$ws.Range("B28").Value = "LYG"
$ws.Range("C28").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D28").Value = 13.5
$ws.Range("A29").Value = "2019Q4"
$ws.Range("B29").Value = "LYG"
$ws.Range("C29").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D29").Value = 13.8
$ws.Range("A30").Value = "2020Q1"
$ws.Range("B30").Value = "LYG"
$ws.Range("C30").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D30").Value = 14.2
$ws.Range("A31").Value = "2020Q2"
$ws.Range("B31").Value = "LYG"
$ws.Range("C31").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D31").Value = 14.6
$ws.Range("A32").Value = "2020Q3"
$ws.Range("B32").Value = "LYG"
$ws.Range("C32").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D32").Value = 15.2
$ws.Range("A33").Value = "2020Q4"
$ws.Range("B33").Value = "LYG"
$ws.Range("C33").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D33").Value = 16.2

# --- Rows 34-46: new unstyled rows, quarters 2021Q1-2024Q1 ---
$ws.Range("A34").Value = "2021Q1"
$ws.Range("B34").Value = "LYG"
$ws.Range("C34").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D34").Value = 16.7
$ws.Range("A35").Value = "2021Q2"
$ws.Range("B35").Value = "LYG"
$ws.Range("C35").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D35").Value = 16.7
$ws.Range("A36").Value = "2021Q3"
$ws.Range("B36").Value = "LYG"
$ws.Range("C36").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D36").Value = 17.2
$ws.Range("A37").Value = "2021Q4"
$ws.Range("B37").Value = "LYG"
$ws.Range("C37").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D37").Value = 16.3
$ws.Range("A38").Value = "2022Q1"
$ws.Range("B38").Value = "LYG"
$ws.Range("C38").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D38").Value = 14.2
$ws.Range("A39").Value = "2022Q2"
$ws.Range("B39").Value = "LYG"
$ws.Range("C39").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D39").Value = 14.7
$ws.Range("A40").Value = "2022Q3"
$ws.Range("B40").Value = "LYG"
$ws.Range("C40").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D40").Value = 15
$ws.Range("A41").Value = "2022Q4"
$ws.Range("B41").Value = "LYG"
$ws.Range("C41").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D41").Value = 14.1
$ws.Range("A42").Value = "2023Q1"
$ws.Range("B42").Value = "LYG"
$ws.Range("C42").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D42").Value = 14.1
$ws.Range("A43").Value = "2023Q2"
$ws.Range("B43").Value = "LYG"
$ws.Range("C43").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D43").Value = 14.2
$ws.Range("A44").Value = "2023Q3"
$ws.Range("B44").Value = "LYG"
$ws.Range("C44").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D44").Value = 14.6
$ws.Range("A45").Value = "2023Q4"
$ws.Range("B45").Value = "LYG"
$ws.Range("C45").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D45").Value = 13.7
$ws.Range("A46").Value = "2024Q1"
$ws.Range("B46").Value = "LYG"
$ws.Range("C46").Value = "LLOYDS BANKING GROUP PLC"
$ws.Range("D46").Value = 13.7

# --- Match final selection state ---
$ws.Activate()
$ws.Range("A47").Select()

Write-Output "Edit complete"